$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$clothing = @{
    2 = 'Trunks,Jumpsuit'
    3 = 'Trunks,Dress'
    4 = 'Jumpsuit,Trunks'
    5 = 'Sweatpants,Dress'
    6 = 'Blazer,Blouse'
    7 = 'Parka,Blouse'
    8 = 'Jumpsuit,Kaftan'
    9 = 'Jumpsuit,Dress'
    10 = 'Tee,Halter'
    11 = 'Blazer,Halter'
    12 = 'Dress,Jumpsuit'
    13 = 'Jumpsuit,Blouse'
    14 = 'Jumpsuit,Blouse'
    15 = 'Blouse,Jumpsuit'
    16 = 'Parka,Jumpsuit'
    17 = 'Caftan,Blazer'
    18 = 'Blouse,Blazer'
    19 = 'Blazer,Caftan'
    20 = 'Jumpsuit,Kaftan'
    21 = 'Jumpsuit,Kaftan'
    22 = 'Jumpsuit,Trunks'
    23 = 'Jumpsuit,Tee'
    24 = 'Jumpsuit,Blouse'
    25 = 'Trunks,Kaftan'
    26 = 'Parka,Jumpsuit'
    27 = 'Jumpsuit,Trunks'
    28 = 'Jumpsuit,Kaftan'
    29 = 'Jumpsuit,Caftan'
    30 = 'Jumpsuit,Blouse'
    31 = 'Blazer,Jumpsuit'
    32 = 'Jumpsuit,Dress'
    33 = 'Jumpsuit,Blouse'
    34 = 'Caftan,Trunks'
    35 = 'Halter,Blazer'
    36 = 'Halter,Blazer'
    37 = 'Halter,Jumpsuit'
    38 = 'Blazer,Blouse'
    39 = 'Trunks,Kaftan'
    40 = 'Parka,Blouse'
    41 = 'Halter,Blazer'
    42 = 'Jumpsuit,Kaftan'
    43 = 'Jumpsuit,Kaftan'
    44 = 'Dress,Kaftan'
    45 = 'Jodhpurs,Trunks'
    46 = 'Parka,Sweatpants'
    47 = 'Jumpsuit,Blouse'
    48 = 'Blouse,Jumpsuit'
    49 = 'Trunks,Dress'
    50 = 'Trunks,Jodhpurs'
    51 = 'Jumpsuit,Caftan'
    52 = 'Jumpsuit,Halter'
    53 = 'Jumpsuit,Kaftan'
    54 = 'Halter,Blazer'
    55 = 'Jumpsuit,Halter'
    56 = 'Trunks,Kaftan'
    57 = 'Trunks,Jumpsuit'
    58 = 'Sweatpants,Jumpsuit'
    59 = 'Jumpsuit,Blouse'
    60 = 'Caftan,Jodhpurs'
    61 = 'Trunks,Blazer'
    62 = 'Jumpsuit,Trunks'
    63 = 'Turtleneck,Jodhpurs'
    64 = 'Trunks,Jumpsuit'
    65 = 'Halter,Jumpsuit'
    66 = 'Halter,Sweatpants'
    67 = 'Caftan,Parka'
    68 = 'Blazer,Blouse'
    69 = 'Parka,Dress'
    70 = 'Parka,Sweatpants'
    71 = 'Parka,Sweatpants'
    72 = 'Dress,Trunks'
    73 = 'Caftan,Jumpsuit'
    74 = 'Jumpsuit,Coverup'
    75 = 'Cutoffs,Jodhpurs'
    76 = 'Jumpsuit,Dress'
    77 = 'Jumpsuit,Parka'
    78 = 'Blouse,Jumpsuit'
    79 = 'Trunks,Blazer'
    80 = 'Parka,Sweatpants'
    81 = 'Sweatpants,Parka'
    82 = 'Caftan,Blouse'
    83 = 'Halter,Blouse'
    84 = 'Halter,Jumpsuit'
    85 = 'Jumpsuit,Halter'
    86 = 'Halter,Caftan'
    87 = 'Jumpsuit,Trunks'
    88 = 'Blouse,Kaftan'
    89 = 'Jumpsuit,Blouse'
    90 = 'Jumpsuit,Blouse'
    91 = 'Jumpsuit,Blouse'
    92 = 'Trunks,Kaftan'
    93 = 'Halter,Jumpsuit'
    94 = 'Jumpsuit,Trunks'
    95 = 'Jumpsuit,Dress'
    96 = 'Trunks,Jumpsuit'
    97 = 'Blouse,Halter'
    98 = 'Blouse,Trunks'
    99 = 'Jumpsuit,Caftan'
    100 = 'Blazer,Trunks'
    101 = 'Jumpsuit,Kaftan'
    102 = 'Halter,Top'
    103 = 'Jumpsuit,Halter'
    104 = 'Halter,Top'
    105 = 'Jumpsuit,Blouse'
    106 = 'Blazer,Top'
    107 = 'Caftan,Trunks'
    108 = 'Trunks,Blouse'
    109 = 'Jumpsuit,Kaftan'
    110 = 'Halter,Jumpsuit'
    111 = 'Blouse,Jumpsuit'
    112 = 'Parka,Blouse'
    113 = 'Jumpsuit,Halter'
    114 = 'Jumpsuit,Kaftan'
    115 = 'Blazer,Top'
    116 = 'Kaftan,Jumpsuit'
}

foreach ($row in $clothing.Keys) {
    $ws.Range("G$row").Value = $clothing[$row]
}
